# Latest data up to Nov 2023
# Update the Poker "Year Figures" sheet, recomputed standings for year 2023
# (rows 180-189) to reflect the latest data through November 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, Person, SRank, Points, Chips(Bonus col E), PointsBonus col F, Winnings, Takehome, (I) delta, Takehome total (K)
$rows = @(
    @{ Row=180; B="Jon";      C=1;  D=35; E=13; F=48; G=117250; H=60;  I=-40; K=357 },
    @{ Row=181; B="Matt";     C=2;  D=41; E=6;  F=47; G=128650; H=70;  I=-30; K=362 },
    @{ Row=182; B="Maisy";    C=3;  D=39; E=6;  F=45; G=127200; H=100; I=20;  K=360 },
    @{ Row=183; B="Anthony";  C=4;  D=38; E=4;  F=42; G=133800; H=170; I=100; K=350 },
    @{ Row=184; B="Pepe";     C=5;  D=29; E=3;  F=32; G=90100;  H=80;  I=10;  K=364 },
    @{ Row=185; B="Richard";  C=6;  D=28; E=2;  F=30; G=93050;  H=60;  I=-30; K=366 },
    @{ Row=186; B="Andy";     C=7;  D=21; E=6;  F=27; G=77700;  H=80;  I=-10; K=349 },
    @{ Row=187; B="Douglas";  C=8;  D=24; E=3;  F=27; G=71400;  H=100; I=50;  K=424 },
    @{ Row=188; B="Prashant"; C=9;  D=23; E=3;  F=26; G=78850;  H=90;  I=20;  K=365 },
    @{ Row=189; B="Mark";     C=10; D=21; E=2;  F=23; G=63200;  H=10;  I=-60; K=361 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value  = $r.B   # B: Person
    $ws.Cells.Item($n, 3).Value  = $r.C   # C: SRank
    $ws.Cells.Item($n, 4).Value  = $r.D   # D: Points
    $ws.Cells.Item($n, 5).Value  = $r.E   # E: Bonus
    $ws.Cells.Item($n, 6).Value  = $r.F   # F: PointsBonus
    $ws.Cells.Item($n, 7).Value  = $r.G   # G: Chips/Winnings
    $ws.Cells.Item($n, 8).Value  = $r.H   # H: Takehome
    $ws.Cells.Item($n, 9).Value  = $r.I   # I: delta
    $ws.Cells.Item($n, 11).Value = $r.K   # K: pers_personid
}
